# Insert a new data row right before the current row 242 (Excel shifts
# rows 242..342 down to 243..343), then populate the newly inserted row
# with the new record. This reproduces the diff, which shows a single
# new row of data added in the middle of the table (row 242), pushing
# every subsequent row down by one and growing the sheet from 342 to
# 343 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 242:342 down to 243:343, leaving a blank row 242 that
# inherits formatting (including the date number format) from the row
# above it, same as Excel's normal "Insert" behavior.
$ws.Rows.Item(242).Insert()

$ws.Range("A242").Value = 9
$ws.Range("B242").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C242").Value = "Metropolitana"
$ws.Range("D242").Value = 44837
$ws.Range("D242").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E242").Value = 13
$ws.Range("F242").Value = 100112021
$ws.Range("G242").Value = "Ají"
$ws.Range("H242").Value = "Inferno"
$ws.Range("I242").Value = "Primera"
$ws.Range("J242").Value = 70
$ws.Range("K242").Value = 28000
$ws.Range("L242").Value = 28000
$ws.Range("M242").Value = 28000
$ws.Range("N242").Value = "`$/caja 10 kilos"
$ws.Range("O242").Value = "Región de Arica y Parinacota"
$ws.Range("P242").Value = 2800
$ws.Range("Q242").Value = 10
$ws.Range("R242").Value = "Hortaliza"
